# Refresh the Rspo1-Lgr4 NATMI LR-pair sheet with newly computed TPM values.
# The "Resolving-Mac" sending-cluster rows have been dropped, and
# "Inflammatory-Mac" now also appears as a target-cluster option, shrinking
# the table from 12 data rows (A1:T13) to 10 data rows (A1:T11). Every
# numeric metric column (G:T) has been recomputed against the new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two now-obsolete data rows so the used range becomes A1:T11.
$ws.Rows("12:13").Delete()

# Overwrite the remaining data rows (A2:T11) with the refreshed values.
$newData = New-Object "object[,]" 10,20
$newData[0,0] = "ECs"
$newData[0,1] = "Rspo1"
$newData[0,2] = "Lgr4"
$newData[0,3] = "ECs"
$newData[0,4] = 1
$newData[0,5] = 0.3333333333333333
$newData[0,6] = 0.03450166666666667
$newData[0,7] = 0.103505
$newData[0,8] = 0.06325845361026364
$newData[0,9] = 0.06325845361026364
$newData[0,10] = 3
$newData[0,11] = 1
$newData[0,12] = 0.5779736666666667
$newData[0,13] = 1.733921
$newData[0,14] = 0.04835019606981441
$newData[0,15] = 0.05356228267519154
$newData[0,16] = 0.01994105478944445
$newData[0,17] = 0.179469493105
$newData[0,18] = 0.003058558635129506
$newData[0,19] = 0.003388267173868432
$newData[1,0] = "ECs"
$newData[1,1] = "Rspo1"
$newData[1,2] = "Lgr4"
$newData[1,3] = "FAPs"
$newData[1,4] = 1
$newData[1,5] = 0.3333333333333333
$newData[1,6] = 0.03450166666666667
$newData[1,7] = 0.103505
$newData[1,8] = 0.06325845361026364
$newData[1,9] = 0.06325845361026364
$newData[1,10] = 3
$newData[1,11] = 1
$newData[1,12] = 7.791016
$newData[1,13] = 23.373048
$newData[1,14] = 0.6517548686181108
$newData[1,15] = 0.7220131735856595
$newData[1,16] = 0.2688030370266667
$newData[1,17] = 2.41922733324
$newData[1,18] = 0.04122900512174223
$newData[1,19] = 0.04567343684726767
$newData[2,0] = "ECs"
$newData[2,1] = "Rspo1"
$newData[2,2] = "Lgr4"
$newData[2,3] = "Inflammatory-Mac"
$newData[2,4] = 1
$newData[2,5] = 0.3333333333333333
$newData[2,6] = 0.03450166666666667
$newData[2,7] = 0.103505
$newData[2,8] = 0.06325845361026364
$newData[2,9] = 0.06325845361026364
$newData[2,10] = 1
$newData[2,11] = 0.3333333333333333
$newData[2,12] = 0.05015166666666667
$newData[2,13] = 0.150455
$newData[2,14] = 0.004195421100317676
$newData[2,15] = 0.004647681895481942
$newData[2,16] = 0.001730316086111111
$newData[2,17] = 0.015572844775
$newData[2,18] = 0.000265395851049967
$newData[2,19] = 0.0002940051695806066
$newData[3,0] = "ECs"
$newData[3,1] = "Rspo1"
$newData[3,2] = "Lgr4"
$newData[3,3] = "MuSCs"
$newData[3,4] = 1
$newData[3,5] = 0.3333333333333333
$newData[3,6] = 0.03450166666666667
$newData[3,7] = 0.103505
$newData[3,8] = 0.06325845361026364
$newData[3,9] = 0.06325845361026364
$newData[3,10] = 2
$newData[3,11] = 1
$newData[3,12] = 3.489664
$newData[3,13] = 6.979328
$newData[3,14] = 0.2919266886169084
$newData[3,15] = 0.215597330685123
$newData[3,16] = 0.1203992241066667
$newData[3,17] = 0.72239534464
$newData[3,18] = 0.01846683088947058
$newData[3,19] = 0.01363835374164152
$newData[4,0] = "ECs"
$newData[4,1] = "Rspo1"
$newData[4,2] = "Lgr4"
$newData[4,3] = "Resolving-Mac"
$newData[4,4] = 1
$newData[4,5] = 0.3333333333333333
$newData[4,6] = 0.03450166666666667
$newData[4,7] = 0.103505
$newData[4,8] = 0.06325845361026364
$newData[4,9] = 0.06325845361026364
$newData[4,10] = 1
$newData[4,11] = 0.3333333333333333
$newData[4,12] = 0.0451
$newData[4,13] = 0.1353
$newData[4,14] = 0.003772825594848836
$newData[4,15] = 0.004179531158543795
$newData[4,16] = 0.001556025166666667
$newData[4,17] = 0.0140042265
$newData[4,18] = 0.0002386631128713604
$newData[4,19] = 0.0002643906779053941
$newData[5,0] = "Inflammatory-Mac"
$newData[5,1] = "Rspo1"
$newData[5,2] = "Lgr4"
$newData[5,3] = "ECs"
$newData[5,4] = 1
$newData[5,5] = 0.3333333333333333
$newData[5,6] = 0.5109063333333334
$newData[5,7] = 1.532719
$newData[5,8] = 0.9367415463897364
$newData[5,9] = 0.9367415463897364
$newData[5,10] = 3
$newData[5,11] = 1
$newData[5,12] = 0.5779736666666667
$newData[5,13] = 1.733921
$newData[5,14] = 0.04835019606981441
$newData[5,15] = 0.05356228267519154
$newData[5,16] = 0.2952904067998889
$newData[5,17] = 2.657613661199
$newData[5,18] = 0.0452916374346849
$newData[5,19] = 0.05017401550132311
$newData[6,0] = "Inflammatory-Mac"
$newData[6,1] = "Rspo1"
$newData[6,2] = "Lgr4"
$newData[6,3] = "FAPs"
$newData[6,4] = 1
$newData[6,5] = 0.3333333333333333
$newData[6,6] = 0.5109063333333334
$newData[6,7] = 1.532719
$newData[6,8] = 0.9367415463897364
$newData[6,9] = 0.9367415463897364
$newData[6,10] = 3
$newData[6,11] = 1
$newData[6,12] = 7.791016
$newData[6,13] = 23.373048
$newData[6,14] = 0.6517548686181108
$newData[6,15] = 0.7220131735856595
$newData[6,16] = 3.980479417501333
$newData[6,17] = 35.824314757512
$newData[6,18] = 0.6105258634963686
$newData[6,19] = 0.6763397367383919
$newData[7,0] = "Inflammatory-Mac"
$newData[7,1] = "Rspo1"
$newData[7,2] = "Lgr4"
$newData[7,3] = "Inflammatory-Mac"
$newData[7,4] = 1
$newData[7,5] = 0.3333333333333333
$newData[7,6] = 0.5109063333333334
$newData[7,7] = 1.532719
$newData[7,8] = 0.9367415463897364
$newData[7,9] = 0.9367415463897364
$newData[7,10] = 1
$newData[7,11] = 0.3333333333333333
$newData[7,12] = 0.05015166666666667
$newData[7,13] = 0.150455
$newData[7,14] = 0.004195421100317676
$newData[7,15] = 0.004647681895481942
$newData[7,16] = 0.02562280412722222
$newData[7,17] = 0.230605237145
$newData[7,18] = 0.00393002524926771
$newData[7,19] = 0.004353676725901335
$newData[8,0] = "Inflammatory-Mac"
$newData[8,1] = "Rspo1"
$newData[8,2] = "Lgr4"
$newData[8,3] = "MuSCs"
$newData[8,4] = 1
$newData[8,5] = 0.3333333333333333
$newData[8,6] = 0.5109063333333334
$newData[8,7] = 1.532719
$newData[8,8] = 0.9367415463897364
$newData[8,9] = 0.9367415463897364
$newData[8,10] = 2
$newData[8,11] = 1
$newData[8,12] = 3.489664
$newData[8,13] = 6.979328
$newData[8,14] = 0.2919266886169084
$newData[8,15] = 0.215597330685123
$newData[8,16] = 1.782891438805333
$newData[8,17] = 10.697348632832
$newData[8,18] = 0.2734598577274379
$newData[8,19] = 0.2019589769434815
$newData[9,0] = "Inflammatory-Mac"
$newData[9,1] = "Rspo1"
$newData[9,2] = "Lgr4"
$newData[9,3] = "Resolving-Mac"
$newData[9,4] = 1
$newData[9,5] = 0.3333333333333333
$newData[9,6] = 0.5109063333333334
$newData[9,7] = 1.532719
$newData[9,8] = 0.9367415463897364
$newData[9,9] = 0.9367415463897364
$newData[9,10] = 1
$newData[9,11] = 0.3333333333333333
$newData[9,12] = 0.0451
$newData[9,13] = 0.1353
$newData[9,14] = 0.003772825594848836
$newData[9,15] = 0.004179531158543795
$newData[9,16] = 0.02304187563333333
$newData[9,17] = 0.2073768807
$newData[9,18] = 0.003534162481977476
$newData[9,19] = 0.003915140480638402

$ws.Range("A2:T11").Value = $newData
